$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 20 de Mayo de 2020 a las 23:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1586638
$ws.Range("C4").Value = 16055
$ws.Range("D4").Value = 361553
$ws.Range("E4").Value = 1130377
$ws.Range("G4").Value = 1175
$ws.Range("H4").Value = 94708

# Row 11 - Alemania
$ws.Range("B11").Value = 178489
$ws.Range("C11").Value = 662
$ws.Range("E11").Value = 13324

# Row 25 - Ecuador
$ws.Range("B25").Value = 34854
$ws.Range("C25").Value = 703
$ws.Range("D25").Value = 3557
$ws.Range("E25").Value = 28409
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = 2888

# Row 118 - Principado de Andorra
$ws.Range("B118").Value = 762
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 639
$ws.Range("E118").Value = 72

# Row 152 - Suazilandia
$ws.Range("B152").Value = 217
$ws.Range("C152").Value = 9
$ws.Range("D152").Value = 97
$ws.Range("E152").Value = 118
